# Update betting odds / correct-score figures on Sheet1 to match the
# refreshed FlashScore export for 2024-11-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 ---
$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 5.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("AJ2").Value = 51
$ws.Range("AL2").Value = 51
$ws.Range("AQ2").Value = 41
$ws.Range("AZ2").Value = 126

# --- Row 4 ---
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5

# --- Row 5 ---
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3

# --- Row 6 ---
$ws.Range("G6").Value = 1.95
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 4.33
$ws.Range("J6").Value = 2.75
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 4.75
$ws.Range("Q6").Value = 2.35
$ws.Range("R6").Value = 1.57
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.38
$ws.Range("X6").Value = 8
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 17
$ws.Range("AA6").Value = 19
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 7
$ws.Range("AE6").Value = 19
$ws.Range("AG6").Value = 9.5
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 15
$ws.Range("AK6").Value = 41
$ws.Range("AL6").Value = 51
$ws.Range("AN6").Value = 3.75
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 26
$ws.Range("AR6").Value = 67
$ws.Range("AS6").Value = 201
$ws.Range("AT6").Value = 2.38
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 23

# --- Row 8 ---
$ws.Range("G8").Value = 2.45
$ws.Range("H8").Value = 2.75
$ws.Range("I8").Value = 3.1
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 5.5
$ws.Range("S8").Value = 1.67
$ws.Range("T8").Value = 2.1
$ws.Range("Z8").Value = 26
$ws.Range("AH8").Value = 13
$ws.Range("AO8").Value = 17
$ws.Range("AT8").Value = 2.1
$ws.Range("AZ8").Value = 67
